# Update column G ("K" = strikeouts) values per regenerated save_data.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" — the K column (strikeouts) is recalculated from the
# source box-score data and rewritten into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1
    3 = 0
    4 = 1
    5 = 1
    6 = 1
    7 = 0
    8 = 1
    9 = 0
    10 = 0
    11 = 0
    12 = 2
    13 = 1
    14 = 0
    15 = 1
    16 = 2
    17 = 2
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 3
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 0
    28 = 0
    29 = 2
    30 = 0
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 0
    38 = 3
    39 = 0
    40 = 2
    42 = 1
    43 = 2
    44 = 1
    45 = 2
    46 = 0
    47 = 0
    48 = 2
    49 = 2
    50 = 3
    51 = 0
    52 = 2
    53 = 3
    54 = 0
    55 = 2
    56 = 2
    57 = 0
    58 = 0
    59 = 0
    61 = 0
    62 = 1
    63 = 0
    64 = 0
    65 = 3
    66 = 1
    67 = 2
    68 = 2
    69 = 2
    70 = 3
    71 = 2
    73 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
